# Ticket 516: import product via part number, rather than display name
#
# Row 2 of the Sales Forecast import template is sample/demo data. Previously
# the "Reference" column (A2) held the product's long display name; it now
# holds the product's part number instead, so the importer can match on
# part number. The other sample values are refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sample data: reference by part number, not display name ---
$ws.Range("A2").Value = "100017-000"
$ws.Range("B2").Value = "Internal"
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = "2021-09-11"

# The old display name was long and wrapped onto multiple lines, forcing a
# taller row; the part number is short, so let the row shrink back down.
$ws.Rows.Item(2).AutoFit()

# --- Column widths: size each column to its own content instead of one
#     uniform width across the whole sheet ---
$ws.Columns.Item(1).ColumnWidth = 20.5
$ws.Columns.Item(2).ColumnWidth = 13
$ws.Columns.Item(4).ColumnWidth = 8.833333333333332
$ws.Columns.Item(5).ColumnWidth = 19.333333333333336
$ws.Columns.Item(6).ColumnWidth = 7.666666666666666

# --- View: scroll back to column A and select the last edited cell ---
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("F2").Select()
